## ================================================================
## 1. Reorganize / create worksheets so the final order/names are:
##    Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
## ================================================================

$wb = $excel.ActiveWorkbook

# --- Add "Player Info" right before the existing "ODI Batting" sheet ---
$battingWsForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingWsForInsert)
$playerInfo.Name = "Player Info"

# --- Add "ODI Batting Extra" right after the existing "ODI Bowling" sheet ---
# (re-fetch ODI Bowling by name since sheet references must be re-resolved
#  after any call that mutates the Worksheets collection)
$bowlingWsForInsert = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowlingWsForInsert)
$battingExtra.Name = "ODI Batting Extra"

## ================================================================
## 2. Style + fill the new "Player Info" sheet
## ================================================================

$playerInfoHeader = $playerInfo.Range("A1:D1")
$playerInfoHeader.Font.Bold = $true
$playerInfoHeader.HorizontalAlignment = -4108   # xlCenter
$playerInfoHeader.VerticalAlignment = -4160     # xlTop
$playerInfoHeader.Borders.LineStyle = 1

## ================================================================
## 3. Style the new "ODI Batting Extra" header row
## ================================================================

$battingExtraHeader = $battingExtra.Range("A1:F1")
$battingExtraHeader.Font.Bold = $true
$battingExtraHeader.HorizontalAlignment = -4108 # xlCenter
$battingExtraHeader.VerticalAlignment = -4160   # xlTop
$battingExtraHeader.Borders.LineStyle = 1

# player info headers
$playerInfo.Cells.Item(1,1).Value = 'ID'
$playerInfo.Cells.Item(1,2).Value = 'NAME'
$playerInfo.Cells.Item(1,3).Value = 'BATTING_HAND'
$playerInfo.Cells.Item(1,4).Value = 'BOWL_STYLE'
# player info data rows
$playerInfo.Cells.Item(2,1).NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = '4816'
$playerInfo.Cells.Item(2,2).NumberFormat = "@"
$playerInfo.Cells.Item(2,2).Value = 'Marnus Labuschagne'
$playerInfo.Cells.Item(2,3).NumberFormat = "@"
$playerInfo.Cells.Item(2,3).Value = 'Right Handed'
$playerInfo.Cells.Item(2,4).NumberFormat = "@"
$playerInfo.Cells.Item(2,4).Value = 'Right Arm Leg Break'

# extra headers
$battingExtra.Cells.Item(1,1).Value = 'MATCH_CODE'
$battingExtra.Cells.Item(1,2).Value = 'BATTING_POSITION'
$battingExtra.Cells.Item(1,3).Value = 'NUM_4'
$battingExtra.Cells.Item(1,4).Value = 'NUM_6'
$battingExtra.Cells.Item(1,5).Value = 'PERCENT_RUNS_OF_TOTAL'
$battingExtra.Cells.Item(1,6).Value = 'MAN_OF_MATCH'

# extra data rows
$battingExtra.Cells.Item(2,1).NumberFormat = "@"
$battingExtra.Cells.Item(2,1).Value = '4435'
$battingExtra.Cells.Item(2,6).NumberFormat = "@"
$battingExtra.Cells.Item(2,6).Value = 'NO'
$battingExtra.Cells.Item(3,1).NumberFormat = "@"
$battingExtra.Cells.Item(3,1).Value = '4436'
$battingExtra.Cells.Item(3,6).NumberFormat = "@"
$battingExtra.Cells.Item(3,6).Value = 'NO'
$battingExtra.Cells.Item(4,1).NumberFormat = "@"
$battingExtra.Cells.Item(4,1).Value = '4437'
$battingExtra.Cells.Item(4,2).Value = 1
$battingExtra.Cells.Item(4,3).NumberFormat = "@"
$battingExtra.Cells.Item(4,3).Value = '1'
$battingExtra.Cells.Item(4,4).NumberFormat = "@"
$battingExtra.Cells.Item(4,4).Value = '0'
$battingExtra.Cells.Item(4,5).NumberFormat = "@"
$battingExtra.Cells.Item(4,5).Value = '2.42%'
$battingExtra.Cells.Item(4,6).NumberFormat = "@"
$battingExtra.Cells.Item(4,6).Value = 'NO'
$battingExtra.Cells.Item(5,1).NumberFormat = "@"
$battingExtra.Cells.Item(5,1).Value = '4564'
$battingExtra.Cells.Item(5,2).Value = 4
$battingExtra.Cells.Item(5,3).NumberFormat = "@"
$battingExtra.Cells.Item(5,3).Value = '1'
$battingExtra.Cells.Item(5,4).NumberFormat = "@"
$battingExtra.Cells.Item(5,4).Value = '0'
$battingExtra.Cells.Item(5,5).NumberFormat = "@"
$battingExtra.Cells.Item(5,5).Value = '7.99%'
$battingExtra.Cells.Item(5,6).NumberFormat = "@"
$battingExtra.Cells.Item(5,6).Value = 'NO'
$battingExtra.Cells.Item(6,1).NumberFormat = "@"
$battingExtra.Cells.Item(6,1).Value = '4565'
$battingExtra.Cells.Item(6,6).NumberFormat = "@"
$battingExtra.Cells.Item(6,6).Value = 'NO'
$battingExtra.Cells.Item(7,1).NumberFormat = "@"
$battingExtra.Cells.Item(7,1).Value = '4567'
$battingExtra.Cells.Item(7,2).Value = 4
$battingExtra.Cells.Item(7,3).NumberFormat = "@"
$battingExtra.Cells.Item(7,3).Value = '1'
$battingExtra.Cells.Item(7,4).NumberFormat = "@"
$battingExtra.Cells.Item(7,4).Value = '0'
$battingExtra.Cells.Item(7,5).NumberFormat = "@"
$battingExtra.Cells.Item(7,5).Value = '1.90%'
$battingExtra.Cells.Item(7,6).NumberFormat = "@"
$battingExtra.Cells.Item(7,6).Value = 'NO'
$battingExtra.Cells.Item(8,1).NumberFormat = "@"
$battingExtra.Cells.Item(8,1).Value = '4594'
$battingExtra.Cells.Item(8,6).NumberFormat = "@"
$battingExtra.Cells.Item(8,6).Value = 'NO'
$battingExtra.Cells.Item(9,1).NumberFormat = "@"
$battingExtra.Cells.Item(9,1).Value = '4597'
$battingExtra.Cells.Item(9,2).Value = 5
$battingExtra.Cells.Item(9,3).NumberFormat = "@"
$battingExtra.Cells.Item(9,3).Value = '0'
$battingExtra.Cells.Item(9,4).NumberFormat = "@"
$battingExtra.Cells.Item(9,4).Value = '0'
$battingExtra.Cells.Item(9,5).NumberFormat = "@"
$battingExtra.Cells.Item(9,5).Value = '9.52%'
$battingExtra.Cells.Item(9,6).NumberFormat = "@"
$battingExtra.Cells.Item(9,6).Value = 'NO'
$battingExtra.Cells.Item(10,1).NumberFormat = "@"
$battingExtra.Cells.Item(10,1).Value = '4600'
$battingExtra.Cells.Item(10,6).NumberFormat = "@"
$battingExtra.Cells.Item(10,6).Value = 'NO'
$battingExtra.Cells.Item(11,1).NumberFormat = "@"
$battingExtra.Cells.Item(11,1).Value = '4601'
$battingExtra.Cells.Item(11,2).Value = 4
$battingExtra.Cells.Item(11,3).NumberFormat = "@"
$battingExtra.Cells.Item(11,3).Value = '0'
$battingExtra.Cells.Item(11,4).NumberFormat = "@"
$battingExtra.Cells.Item(11,4).Value = '0'
$battingExtra.Cells.Item(11,5).NumberFormat = "@"
$battingExtra.Cells.Item(11,5).Value = '5.51%'
$battingExtra.Cells.Item(11,6).NumberFormat = "@"
$battingExtra.Cells.Item(11,6).Value = 'NO'
$battingExtra.Cells.Item(12,1).NumberFormat = "@"
$battingExtra.Cells.Item(12,1).Value = '4603'
$battingExtra.Cells.Item(12,2).Value = 5
$battingExtra.Cells.Item(12,3).NumberFormat = "@"
$battingExtra.Cells.Item(12,3).Value = '2'
$battingExtra.Cells.Item(12,4).NumberFormat = "@"
$battingExtra.Cells.Item(12,4).Value = '0'
$battingExtra.Cells.Item(12,5).NumberFormat = "@"
$battingExtra.Cells.Item(12,5).Value = '18.90%'
$battingExtra.Cells.Item(12,6).NumberFormat = "@"
$battingExtra.Cells.Item(12,6).Value = 'NO'
$battingExtra.Cells.Item(13,1).NumberFormat = "@"
$battingExtra.Cells.Item(13,1).Value = '4647'
$battingExtra.Cells.Item(13,6).NumberFormat = "@"
$battingExtra.Cells.Item(13,6).Value = 'NO'
$battingExtra.Cells.Item(14,1).NumberFormat = "@"
$battingExtra.Cells.Item(14,1).Value = '4648'
$battingExtra.Cells.Item(14,2).Value = 4
$battingExtra.Cells.Item(14,3).NumberFormat = "@"
$battingExtra.Cells.Item(14,3).Value = '0'
$battingExtra.Cells.Item(14,4).NumberFormat = "@"
$battingExtra.Cells.Item(14,4).Value = '0'
$battingExtra.Cells.Item(14,5).NumberFormat = "@"
$battingExtra.Cells.Item(14,5).Value = '2.56%'
$battingExtra.Cells.Item(14,6).NumberFormat = "@"
$battingExtra.Cells.Item(14,6).Value = 'NO'
$battingExtra.Cells.Item(15,1).NumberFormat = "@"
$battingExtra.Cells.Item(15,1).Value = '4649'
$battingExtra.Cells.Item(15,2).Value = 4
$battingExtra.Cells.Item(15,3).NumberFormat = "@"
$battingExtra.Cells.Item(15,3).Value = '2'
$battingExtra.Cells.Item(15,4).NumberFormat = "@"
$battingExtra.Cells.Item(15,4).Value = '0'
$battingExtra.Cells.Item(15,5).NumberFormat = "@"
$battingExtra.Cells.Item(15,5).Value = '19.48%'
$battingExtra.Cells.Item(15,6).NumberFormat = "@"
$battingExtra.Cells.Item(15,6).Value = 'NO'
$battingExtra.Cells.Item(16,1).NumberFormat = "@"
$battingExtra.Cells.Item(16,1).Value = '4660'
$battingExtra.Cells.Item(16,6).NumberFormat = "@"
$battingExtra.Cells.Item(16,6).Value = 'NO'
$battingExtra.Cells.Item(17,1).NumberFormat = "@"
$battingExtra.Cells.Item(17,1).Value = '4663'
$battingExtra.Cells.Item(17,6).NumberFormat = "@"
$battingExtra.Cells.Item(17,6).Value = 'NO'
$battingExtra.Cells.Item(18,1).NumberFormat = "@"
$battingExtra.Cells.Item(18,1).Value = '4666'
$battingExtra.Cells.Item(18,6).NumberFormat = "@"
$battingExtra.Cells.Item(18,6).Value = 'NO'
$battingExtra.Cells.Item(19,1).NumberFormat = "@"
$battingExtra.Cells.Item(19,1).Value = '4725'
$battingExtra.Cells.Item(19,2).Value = 4
$battingExtra.Cells.Item(19,3).NumberFormat = "@"
$battingExtra.Cells.Item(19,3).Value = '1'
$battingExtra.Cells.Item(19,4).NumberFormat = "@"
$battingExtra.Cells.Item(19,4).Value = '0'
$battingExtra.Cells.Item(19,5).NumberFormat = "@"
$battingExtra.Cells.Item(19,5).Value = '7.98%'
$battingExtra.Cells.Item(19,6).NumberFormat = "@"
$battingExtra.Cells.Item(19,6).Value = 'NO'
$battingExtra.Cells.Item(20,1).NumberFormat = "@"
$battingExtra.Cells.Item(20,1).Value = '4728'
$battingExtra.Cells.Item(20,2).Value = 4
$battingExtra.Cells.Item(20,6).NumberFormat = "@"
$battingExtra.Cells.Item(20,6).Value = 'NO'
$battingExtra.Cells.Item(21,1).NumberFormat = "@"
$battingExtra.Cells.Item(21,1).Value = '4732'
$battingExtra.Cells.Item(21,2).Value = 5
$battingExtra.Cells.Item(21,3).NumberFormat = "@"
$battingExtra.Cells.Item(21,3).Value = '1'
$battingExtra.Cells.Item(21,4).NumberFormat = "@"
$battingExtra.Cells.Item(21,4).Value = '1'
$battingExtra.Cells.Item(21,5).NumberFormat = "@"
$battingExtra.Cells.Item(21,5).Value = '10.41%'
$battingExtra.Cells.Item(21,6).NumberFormat = "@"
$battingExtra.Cells.Item(21,6).Value = 'NO'

## ================================================================
## 4. Fix up "ODI Batting" sheet:
##    - rename header D1 MATCH_CARD_LINK -> MATCH_CODE
##    - shrink the MATCH_CARD_LINK urls down to the bare match code
##    - drop the two blank INNING_NUMBER cells (B2 / B30) entirely
## ================================================================

$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $odiBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $odiBatting.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '^.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
    }
}

$odiBatting.Range("B2").Value = $null
$odiBatting.Range("B30").Value = $null

## ================================================================
## 5. Fix up "ODI Bowling" sheet:
##    - rename header B1 MATCH_CARD_LINK -> MATCH_CODE
##    - shrink the MATCH_CARD_LINK urls down to the bare match code
## ================================================================

$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $odiBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $odiBowling.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '^.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
    }
}
